# Updates cryptos list values (Price and Volume(1h) columns) to match latest data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.495.10"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "2.419.37"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.48"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.24"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.11"
$ws.Range("E10").Value = "  +2.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.47"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").Value = "2.412.95"
$ws.Range("E16").Value = "  +3.44%  "

$ws.Range("E17").Value = "  +2.96%  "

$ws.Range("D18").Value = "43.456.21"
$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.12"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.64"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.96"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.44"
$ws.Range("E29").Value = "  +3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.22"
$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +19.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.47"
$ws.Range("E32").Value = "  +6.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.13"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0752"
$ws.Range("E35").Value = "  +3.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "132.26"
$ws.Range("E36").Value = "  +28.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("E38").Value = "  +4.15%  "

$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.04"
$ws.Range("E42").Value = "  -7.19%  "

$ws.Range("D43").Value = "1.944.19"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.30"
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").Value = "2.648.54"
$ws.Range("E48").Value = "  +2.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.57"
$ws.Range("E49").Value = "  +3.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.63"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.35"
